$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.822.00'
$ws.Range("E2").Value = '  -2.07%  '
$ws.Range("D3").Value = '2.409.98'
$ws.Range("E3").Value = '  -1.33%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '''568.51'
$ws.Range("E5").Value = '  -2.55%  '
$ws.Range("D6").Value = '''138.43'
$ws.Range("E6").Value = '  -3.02%  '
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("E8").Value = '  -0.65%  '
$ws.Range("D9").Value = '2.392.47'
$ws.Range("E10").Value = '  -2.04%  '
$ws.Range("E11").Value = '  -0.58%  '
$ws.Range("D12").Value = '''5.05'
$ws.Range("E12").Value = '  -3.01%  '
$ws.Range("D13").Value = '''0.336'
$ws.Range("E13").Value = '  -1.94%  '
$ws.Range("D14").Value = '''25.88'
$ws.Range("E14").Value = '  -2.05%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '2.856.70'
$ws.Range("E15").Value = '  -1.06%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '''0.0000170'
$ws.Range("E16").Value = '  -3.45%  '
$ws.Range("D17").Value = '60.778.50'
$ws.Range("E17").Value = '  -1.98%  '
$ws.Range("D18").Value = '2.400.87'
$ws.Range("E18").Value = '  -1.40%  '
$ws.Range("D19").Value = '''7.77'
$ws.Range("E19").Value = '  +8.82%  '
$ws.Range("D20").Value = '''10.58'
$ws.Range("E20").Value = '  -1.62%  '
$ws.Range("D21").Value = '''321.93'
$ws.Range("E21").Value = '  -1.23%  '
$ws.Range("D22").Value = '''4.02'
$ws.Range("E22").Value = '  -1.82%  '
$ws.Range("E23").Value = '  +1.87%  '
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("D25").Value = '''1.82'
$ws.Range("E25").Value = '  -4.62%  '
$ws.Range("D26").Value = '''64.81'
$ws.Range("D27").Value = '''576.74'
$ws.Range("E27").Value = '  -3.50%  '
$ws.Range("D28").Value = '''8.27'
$ws.Range("E28").Value = '  -9.26%  '
$ws.Range("E29").Value = '  -1.16%  '
$ws.Range("D30").Value = '0.0₃0922'
$ws.Range("E30").Value = '  -4.33%  '
$ws.Range("D31").Value = '''7.82'
$ws.Range("E31").Value = '  -1.78%  '
$ws.Range("E32").Value = '  -5.04%  '
$ws.Range("E33").Value = '  -3.72%  '
$ws.Range("E34").Value = '  -2.85%  '
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").Value = '''151.58'
$ws.Range("E36").Value = '  -1.04%  '
$ws.Range("E37").Value = '  -1.93%  '
$ws.Range("B38").Value = 'NEARProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").Value = '''4.57'
$ws.Range("E38").Value = '  -6.22%  '
$ws.Range("D39").Value = '''0.365'
$ws.Range("E39").Value = '  -2.35%  '
$ws.Range("D40").Value = '''18.15'
$ws.Range("E40").Value = '  -1.28%  '
$ws.Range("D41").Value = '''5.09'
$ws.Range("E41").Value = '  -3.23%  '
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("D43").Value = '''1.65'
$ws.Range("E43").Value = '  -3.04%  '
$ws.Range("D44").Value = '''41.12'
$ws.Range("E44").Value = '  -4.71%  '
$ws.Range("D45").Value = '''2.30'
$ws.Range("E45").Value = '  -7.65%  '
$ws.Range("D46").Value = '''141.59'
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").Value = '0.0₆0262'
$ws.Range("E47").Value = '  -1.21%  '
$ws.Range("E48").Value = '  -3.87%  '
$ws.Range("D49").Value = '''0.583'
$ws.Range("E49").Value = '  -2.82%  '
$ws.Range("B50").Value = 'Hedera'
$ws.Range("C50").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D50").Value = '''0.0500'
$ws.Range("E50").Value = '  -3.52%  '
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").Value = '''19.25'
$ws.Range("E51").Value = '  -2.67%  '
